$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (S001001 / 狐狸): add the long message body in the "content" column (E)
$ws.Range("E4").Value = "未来的大明星，你好！想要在贝壳市找到自己的梦想，钱是不可或缺的。我现在向你借10万，等你一年后发迹了记得还我，祝你好运。"

# Row 5: finish the S001002 entry - fill in sender (C) and the long message body (E)
$ws.Range("C5").Value = "贝壳银行"
$ws.Range("E5").Value = "【贝壳银行】尊敬的客户您好！您现在在我行有贷款，请尽快在规定期限内还清，以免对您造成不便。贝壳银行，真诚为您服务！"

# The long message text wraps within the cell, so grow the row height to fit it
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 57.6

# Move the active selection to match the edited workbook's final cursor position
$ws.Range("F5").Select()
